# productImportExample2.xlsx — commit: "added product_option_id and had
# temporary quick fix for library product price -1"
#
# Column B holds product_option_id. Bump each value by 43 so the imported
# rows continue the id sequence from the previous batch (226-234 -> 269-277).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newProductOptionIds = @(269, 270, 271, 272, 273, 274, 275, 276, 277)

for ($i = 0; $i -lt $newProductOptionIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newProductOptionIds[$i]
}

# Reflect the saved view state: selection moved to the product_option_id
# column (B2:B10) and the window scrolled right so column O is leftmost.
$ws.Activate()
$ws.Range("B2:B10").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 15
